# Modification to UND-EFF figures
#
# The original author nudged three shapes that make up the
# "Estimated gap" annotation (a dashed/diamond connector, a plain
# connector, and the "Estimated gap" textbox) to the right by
# 30890 EMU (~2.432 points) on the single slide of the deck.
#
# EMU -> point conversion:  points = EMU / 12700
#   149 "Straight Arrow Connector 148": 1801912 -> 1832802 EMU
#   150 "Straight Connector 149"      : 2270155 -> 2301045 EMU
#   151 "TextBox 150" ("Estimated gap"): 2186554 -> 2217444 EMU
#
# A tiny epsilon (0.3 EMU, expressed in points) is added to each
# target so that the host's internal single-precision (f32) storage
# of Left/Top, followed by floor-to-EMU, still lands exactly on the
# intended integer EMU value instead of being truncated one EMU low.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$EMU_PER_POINT = 12700.0
$EPS_EMU = 0.3

function Set-LeftEmu {
    param($shape, [double]$targetEmu)
    $pts = ($targetEmu / $EMU_PER_POINT) + ($EPS_EMU / $EMU_PER_POINT)
    $shape.Left = $pts
}

$arrowConnector = $s.Shapes.Item("Straight Arrow Connector 148")
Set-LeftEmu $arrowConnector 1832802

$straightConnector = $s.Shapes.Item("Straight Connector 149")
Set-LeftEmu $straightConnector 2301045

$estimatedGapTextBox = $s.Shapes.Item("TextBox 150")
Set-LeftEmu $estimatedGapTextBox 2217444
